# 2023-08-20 update next cosmetic
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Next_Cosmetic")

# G4: shop candidate id 822 -> 1014
$ws.Range("G4").Value = 1014

# Row 6: cosmetic title/notes changes from "Lelouch\nCode Geass" to "Sweat Beast"
$ws.Range("D6").Value = "Sweat Beast"

# G6: full head image link added
$ws.Range("G6").Value = "https://cdn.discordapp.com/attachments/699111007649398865/1132600250603937842/beast_head_2.png"

# G7: hat image link removed
$ws.Range("G7").Value = ""

# G9: neck image link replaced
$ws.Range("G9").Value = "https://cdn.discordapp.com/attachments/699111007649398865/1132600250847215656/beast_skin_1.png"

# Row 12: creator info updated
$ws.Range("C12").Value = "clown_noes666"
$ws.Range("D12").Value = 469660616
$ws.Range("E12").Value = "115 days"

# Row 14: discord alias updated
$ws.Range("C14").Value = "clown noes666"

# Row 15: discord alias added
$ws.Range("C15").Value = "Clown_noes666"

# L16: shop item row id 822 -> 1014
$ws.Range("L16").Value = 1014
